# Adding new test data (client name / client email) to Sheet1, and
# recording a phonetic-guide font setting that Excel persisted for the
# sheet (this mirrors the extra "size 8 Calibri" font entry that shows
# up in styles.xml together with the phoneticPr on the worksheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update the test data values -----------------------------------
$ws.Range("A2").Value = "jenkins"
$ws.Range("B2").Value = "hendi39@qh.com4"

# --- register the (8pt) phonetic guide font used by the sheet -------
# This causes a new font to be appended to the workbook's font table
# (matching the <phoneticPr fontId="2".../> that Excel writes), without
# altering the existing cell formatting/styles of A2.
$phon = $ws.Range("A2").Phonetics
$phon.Font.Size = 8
$ws.Range("A2").ClearFormats()

# --- update the active selection ------------------------------------
$ws.Range("B6").Select()
